$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# lc.155 and lc.641 passed: Sheet2!B5 and Sheet2!B6 swap values
# (B5 was "25,239", B6 was "155,641" -> now B5="155,641", B6="25,239")
$ws2.Range("B5").Value = "155,641"
$ws2.Range("B6").Value = "25,239"

# Update selections to match the edited workbook state
$ws1.Range("D7").Select()
$ws2.Range("D15").Select()
